$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 90, shifting the existing rows 90-97 down to 91-98.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new weekly data point.
$ws.Cells.Item(90, 1).Value = 10
$ws.Cells.Item(90, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(90, 3).Value = "La Araucanía"
$ws.Cells.Item(90, 4).Value = 44858
$ws.Cells.Item(90, 5).Value = 9
$ws.Cells.Item(90, 6).Value = 100112022
$ws.Cells.Item(90, 7).Value = "Arveja Verde"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 215
$ws.Cells.Item(90, 11).Value = 18000
$ws.Cells.Item(90, 12).Value = 20000
$ws.Cells.Item(90, 13).Value = 19163
$ws.Cells.Item(90, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(90, 15).Value = "Región Metropolitana"
$ws.Cells.Item(90, 16).Value = 767
$ws.Cells.Item(90, 17).Value = 25
$ws.Cells.Item(90, 18).Value = "Hortaliza"
